# surat_keterangan_asal_usul.docx edit
#
# Commit: "add alternate fields in surat (penduduk selector type only)"
#
# Substantive content changes:
#   1. The header word "KANTOR" was previously split into two runs ("KAN" /
#      "TOR") straddling a leftover "_GoBack" bookmark; merge it back into
#      a single "KANTOR" run (Find/Replace across the span removes the
#      bookmark and merges the runs automatically).
#   2. Rename the "dusun" placeholder to "nama_dusun" in the three address
#      blocks that reference it: the subject (penduduk), the father
#      (form.ayah) and the mother (form.ibu).

$d = $word.ActiveDocument

# 1) "KAN" + bookmark + "TOR" -> single run "KANTOR"
$d.Content.Find.Execute("KANTOR", $false, $false, $false, $false, $false, `
    $true, 1, $false, "KANTOR", 2) | Out-Null

# 2) penduduk.dusun -> penduduk.nama_dusun
$d.Content.Find.Execute("{penduduk.dusun}", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{penduduk.nama_dusun}", 2) | Out-Null

# 3) form.ayah.dusun -> form.ayah.nama_dusun
$d.Content.Find.Execute("{form.ayah.dusun}", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{form.ayah.nama_dusun}", 2) | Out-Null

# 4) form.ibu.dusun -> form.ibu.nama_dusun
$d.Content.Find.Execute("{form.ibu.dusun}", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{form.ibu.nama_dusun}", 2) | Out-Null
